$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2: currency "ДЕН" retyped as "ДЕННННН" (font explicitly reset in the process,
# producing a new style override on top of the "Normal 2" cell style)
$ws.Range("C2").Value = "ДЕННННН"
$ws.Range("C2").Font.Name = "Calibri"

# F4: interest-rate number 0.01 overwritten with stray text "дсадсадса"
# (also picks up the same explicit-font style override)
$ws.Range("F4").Value = "дсадсадса"
$ws.Range("F4").Font.Name = "Calibri"

# A handful of other "ДЕН" currency cells were retyped the same way, but kept
# their original cell style (no font override on these)
$ws.Range("C6").Value = "ДЕННННН"
$ws.Range("C10").Value = "ДЕННННН"
$ws.Range("C14").Value = "ДЕННННН"
$ws.Range("C19").Value = "ДЕННННН"

# New explicit column widths for F and G
$ws.Columns.Item(6).ColumnWidth = 57.7109375
$ws.Columns.Item(7).ColumnWidth = 50.28515625

# Final cursor/selection position left on F7
$null = $ws.Range("F7").Select()
